# The diary entry paragraph that used to hold the "weather" text is split
# into two paragraphs:
#   1. A new paragraph (inserted before) that keeps the original "weather"
#      sentence.
#   2. The original paragraph, whose text is replaced with a new, short
#      entry: "使用git创建分支".
#
# Word exposes no scriptable "new paragraph mark inherited text" operation,
# so we reproduce the edit the way a person typing into Word would: place
# the caret at the very start of the target paragraph, type the duplicate
# sentence followed by a paragraph break (which pushes the existing
# paragraph, with its bookmarks intact, one slot down), then find/replace
# the text of that now-shifted paragraph.

$d = $word.ActiveDocument

$weatherText = "天气很热，空调还坏了，又没好好午休，只睡了不到20分钟，对了，这是个dev分支，嗯，是个分支，真的是个分支"
$newText = "使用git创建分支"

# Locate the paragraph that currently holds the long "weather" sentence
# (Range.Text includes the trailing paragraph mark, so compare with
# StartsWith rather than an exact match).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith($weatherText)) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the target paragraph containing the weather text."
}

$target = $d.Paragraphs($targetIndex)

# Collapse the selection to the very beginning of that paragraph and type
# a duplicate of the sentence plus a paragraph break, which inserts a new
# paragraph before it (carrying the same indent/formatting) while leaving
# the original paragraph - and its bookmarks - intact just after it.
$insertionPoint = $d.Range($target.Range.Start, $target.Range.Start)
$insertionPoint.Select()
$word.Selection.TypeText($weatherText)
$word.Selection.TypeText([char]13)

# The original paragraph has now shifted down by one slot; replace its
# text (still the weather sentence) with the new, short entry.
$d2 = $word.ActiveDocument
$origIndex = $targetIndex + 1
$orig = $d2.Paragraphs($origIndex)
$orig.Range.Find.Execute($weatherText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newText, 2)
